$wb = $excel.ActiveWorkbook

# --- Sheet: Correlation matrix ---
$ws1 = $wb.Worksheets.Item("Correlation matrix")

$row2 = @(-0.0578566758318202,0.01572563675836635,0.02217052873220555,-0.03112808543034129,0.009506676583047996,-0.04376408419878811,-0.0133642952007395,0.06926441168555347,0.06041161340274066,0.02852866041351129,-0.06671255020778472,-0.07939371386395846,-0.0136233376951211,-0.02091772575958519,-0.07196674847465961,-0.09351063661446735)
$row3 = @(-0.02796824689522253,0.04783872147277625,0.03070688879037094,-0.01978241853564501,0.0375192878561945,-0.05734591371064292,0.009550133086173319,0.006139371269682767,0.03069759915506819,0.006823753064526763,-0.0341084435056314,-0.0341187653226344,-0.004104351229156675,0.004775875553192661,-0.06006138954307796,-0.03888927236171197)
$row4 = @(-0.02796824689522253,0.04783872147277625,0.03070688879037094,-0.01978241853564501,0.0375192878561945,-0.05734591371064292,0.009550133086173319,0.006139371269682767,0.03069759915506819,0.006823753064526763,-0.0341084435056314,-0.0341187653226344,-0.004104351229156675,0.004775875553192661,-0.06006138954307796,-0.03888927236171197)

for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws1.Cells.Item(2, 2 + $i).Value = $row2[$i]
    $ws1.Cells.Item(3, 2 + $i).Value = $row3[$i]
    $ws1.Cells.Item(4, 2 + $i).Value = $row4[$i]
}

# --- Sheet: Toggles no del ---
$ws2 = $wb.Worksheets.Item("Toggles no del")
$ws2.Range("B2").Value = 310
$ws2.Range("B3").Value = 248
$ws2.Range("B4").Value = 17

# --- Sheet: Toggles del ---
$ws3 = $wb.Worksheets.Item("Toggles del")
$ws3.Range("B2").Value = 263
$ws3.Range("B3").Value = 218
$ws3.Range("B4").Value = 35
$ws3.Range("B5").Value = 59

# --- Sheet: Toggles input del ---
$ws4 = $wb.Worksheets.Item("Toggles input del")
$ws4.Range("B2").Value = 263
$ws4.Range("B3").Value = 218
$ws4.Range("B4").Value = 35
$ws4.Range("B5").Value = 59
